$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Objetivos:" row (row 10) content: was a leftover/placeholder
#     value, now the real objectives paragraph ---
$objetivosTexto = "Possibilitar aos alunos a realização de trabalho de síntese e integração dos conhecimentos adquiridos ao longo do curso, conforme projeto aprovado na disciplina de Trabalho de Conclusão do Curso I."
$ws.Range("B10").Value = $objetivosTexto
$ws.Range("C10").Value = $objetivosTexto

# --- Insert a new row at 13 (pushes "Programa resumido:" and everything
#     below it down by one row, carrying row heights along) so the
#     "Docentes responsáveis:" label (already on row 12) gets its value
#     row right underneath it ---
$ws.Rows.Item(13).Insert()

# The insert carries row 12's (label-only, column-A-styled) formatting
# onto the new row; row 13 has no "A" label though, so drop that stray
# formatted-but-empty cell.
$ws.Range("A13").Clear()

# New row 13 holds the professor name that used to live on row 10.
# B13 needs column B's normal wrap/top-align style pulled in explicitly
# before writing its value (C13 already inherited the right style).
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$professor = "1285870 - Marcos Villela Barcza"
$ws.Range("B13").Value = $professor
$ws.Range("C13").Value = $professor

# Row 14 = "Programa resumido:" -> real short-syllabus text (was "Semestral")
$programaResumido = "Desenvolvimento do trabalho de conclusão de curso, sob orientação de um professor orientador, o qual deve constituir-se num projeto de engenharia química."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 16 = "Programa:" -> real syllabus text (was a stray date value)
$programa = "Elaboração de uma monografia de conclusão de curso que apresente: (1) o tema e sua importância, (2) os objetivos, (3) a revisão bibliográfica, (4) a metodologia científica (5) o desenvolvimento do projeto, (6) a análise e discussão dos resultados, (7) as conclusões e (8) referências bibliográficas."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Row 19 = "Método:" -> method text (was a stray professor-name value)
$metodo = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final, conforme norma do Departamento de Engenharia Química."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 = "Critério:" -> evaluation-board criterion text
$criterio = "Avaliação da monografia perante uma banca examinadora composta por 3 (três) membros, obrigatoriamente docentes da Escola de Engenharia de Lorena (EEL)."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 = "Norma de recuperação:" -> recovery-rule text
$recuperacao = "Reapresentação da monografia, preferencialmente para a mesma banca, com as modificações sugeridas para uma nova avaliação."
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao

# Row 22 = "Bibliografia:" -> bibliography text (was empty before)
$bibliografia = "Recomendada pelo orientador"
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
